$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.669.66'
$ws.Range('E2').Value = '  -0.74%  '
$ws.Range('D3').Value = '1.583.31'
$ws.Range('E3').Value = '  -3.18%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '206.89'
$ws.Range('E5').Value = '  -2.29%  '
$ws.Range('E6').Value = '  -3.12%  '
$ws.Range('E7').Value = '  +0.22%  '
$ws.Range('E8').Value = '  -4.68%  '
$ws.Range('E9').Value = '  -1.37%  '
$ws.Range('E10').Value = '  -3.39%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0868'
$ws.Range('E11').Value = '  -1.81%  '
$ws.Range('D12').Value = '1.808.43'
$ws.Range('E12').Value = '  -3.14%  '
$ws.Range('D13').Value = '1.581.57'
$ws.Range('E13').Value = '  -3.36%  '
$ws.Range('E14').Value = '  -4.04%  '
$ws.Range('E15').Value = '  -5.41%  '
$ws.Range('D16').Value = '27.641.55'
$ws.Range('E16').Value = '  -0.92%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.08'
$ws.Range('E17').Value = '  -3.39%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '217.77'
$ws.Range('E18').Value = '  -4.95%  '
$ws.Range('E19').Value = '  -3.61%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.32'
$ws.Range('E20').Value = '  -4.84%  '
$ws.Range('E21').Value = '  +0.20%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.14'
$ws.Range('E22').Value = '  -4.50%  '
$ws.Range('E23').Value = '  -5.04%  '
$ws.Range('E24').Value = '  -5.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '153.83'
$ws.Range('E25').Value = '  -0.96%  '
$ws.Range('E27').Value = '  -2.59%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.08'
$ws.Range('E28').Value = '  -2.97%  '
$ws.Range('E29').Value = '  -4.23%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.15'
$ws.Range('E30').Value = '  -2.43%  '
$ws.Range('E31').Value = '  -3.53%  '
$ws.Range('E32').Value = '  -5.26%  '
$ws.Range('D33').Value = '1.376.35'
$ws.Range('E33').Value = '  -1.32%  '
$ws.Range('E35').Value = '  -5.33%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.969'
$ws.Range('E36').Value = '  -4.85%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0165'
$ws.Range('E38').Value = '  -3.05%  '
$ws.Range('E39').Value = '  -3.65%  '
$ws.Range('E40').Value = '  -3.89%  '
$ws.Range('E41').Value = '  +0.19%  '
$ws.Range('E42').Value = '  -3.39%  '
$ws.Range('E43').Value = '  -2.78%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.17'
$ws.Range('E44').Value = '  +1.38%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '63.61'
$ws.Range('E45').Value = '  -3.51%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.21'
$ws.Range('E46').Value = '  -4.32%  '
$ws.Range('D47').Value = '1.719.48'
$ws.Range('E47').Value = '  -3.13%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '87.67'
$ws.Range('E48').Value = '  -1.11%  '
$ws.Range('D49').Value = '0.0₆0100'
$ws.Range('E49').Value = '  -2.56%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0974'
$ws.Range('E50').Value = '  -5.00%  '
$ws.Range('E51').Value = '  -1.58%  '
